# Update countries & provincias Spain
# Refresh the COVID-19 country table on sheet "Pais":
#  - bump the "last updated" timestamp in A1
#  - refresh the daily counters (columns B..H) for the countries whose
#    figures moved between the two snapshots
#  - a handful of countries swapped rank (their row's data changed but the
#    two adjacent rows also swap which country name they show), so the A
#    column for those rows is rewritten too

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- timestamp ---------------------------------------------------------
$ws.Range("A1").Value = 'Datos actualizados a 20 de Septiembre de 2020 a las 17:24'

# --- Estados Unidos (row 4) ---------------------------------------------
$ws.Range("B4").Value = 6973062
$ws.Range("C4").Value = 5659
$ws.Range("D4").Value = 4224307
$ws.Range("E4").Value = 2544872
$ws.Range("G4").Value = 59
$ws.Range("H4").Value = 203883

# --- India (row 5) -------------------------------------------------------
$ws.Range("B5").Value = 5417274
$ws.Range("C5").Value = 19044
$ws.Range("D5").Value = 4313402
$ws.Range("E5").Value = 1016963
$ws.Range("G5").Value = 135
$ws.Range("H5").Value = 86909

# --- Reino Unido (row 17) ------------------------------------------------
$ws.Range("B17").Value = 394257
$ws.Range("C17").Value = 3899
$ws.Range("G17").Value = 18
$ws.Range("H17").Value = 41777

# --- Italia (row 23) ------------------------------------------------------
$ws.Range("B23").Value = 298156
$ws.Range("C23").Value = 1587
$ws.Range("D23").Value = 218351
$ws.Range("E23").Value = 44098
$ws.Range("G23").Value = 15
$ws.Range("H23").Value = 35707

# --- Canada (row 29) -------------------------------------------------------
$ws.Range("B29").Value = 143141
$ws.Range("C29").Value = 367
$ws.Range("D29").Value = 124378
$ws.Range("E29").Value = 9551
$ws.Range("G29").Value = 1
$ws.Range("H29").Value = 9212

# --- Singapur (row 57) -----------------------------------------------------
$ws.Range("D57").Value = 57181
$ws.Range("E57").Value = 368

# --- Libia / El Salvador swap ranks (rows 76 & 77) --------------------------
$ws.Range("A76").Value = 'Libia'
$ws.Range("B76").Value = 27949
$ws.Range("C76").Value = 715
$ws.Range("D76").Value = 15068
$ws.Range("E76").Value = 12437
$ws.Range("G76").Value = 8
$ws.Range("H76").Value = 444

$ws.Range("A77").Value = 'El Salvador'
$ws.Range("B77").Value = 27553
$ws.Range("C77").Value = 125
$ws.Range("D77").Value = 21561
$ws.Range("E77").Value = 5181
$ws.Range("G77").Value = 3
$ws.Range("H77").Value = 811

# --- Albania (row 94) --------------------------------------------------------
$ws.Range("B94").Value = 12385
$ws.Range("C94").Value = 159
$ws.Range("D94").Value = 6940
$ws.Range("E94").Value = 5083
$ws.Range("G94").Value = 4
$ws.Range("H94").Value = 362

# --- Namibia (row 96) ----------------------------------------------------------
$ws.Range("B96").Value = 10377
$ws.Range("C96").Value = 85
$ws.Range("D96").Value = 8033
$ws.Range("E96").Value = 2232
$ws.Range("G96").Value = 1
$ws.Range("H96").Value = 112

# --- Jordania / Jamaica / Surinam / Ruanda re-rank (rows 124-127) -------------
$ws.Range("A124").Value = 'Jordania'
$ws.Range("B124").Value = 4779
$ws.Range("C124").Value = 239
$ws.Range("D124").Value = 2844
$ws.Range("E124").Value = 1905
$ws.Range("H124").Value = 30

$ws.Range("A125").Value = 'Jamaica'
$ws.Range("B125").Value = 4758
$ws.Range("D125").Value = 1327
$ws.Range("E125").Value = 3371
$ws.Range("H125").Value = 60

$ws.Range("A126").Value = 'Surinam'
$ws.Range("B126").Value = 4709
$ws.Range("D126").Value = 4383
$ws.Range("E126").Value = 229
$ws.Range("H126").Value = 97

$ws.Range("A127").Value = 'Ruanda'
$ws.Range("B127").Value = 4689
$ws.Range("D127").Value = 2910
$ws.Range("E127").Value = 1753
$ws.Range("H127").Value = 26

# --- Sri Lanka (row 140) -----------------------------------------------------
$ws.Range("B140").Value = 3287
$ws.Range("C140").Value = 4
$ws.Range("E140").Value = 186

# --- Santa Lucia / Timor Oriental swap ranks (rows 204 & 205) ----------------
$ws.Range("A204").Value = 'Santa Lucia'
$ws.Range("A205").Value = 'Timor Oriental'

# --- Montserrat / Islas Malvinas swap ranks (rows 214 & 215) -----------------
$ws.Range("A214").Value = 'Montserrat'
$ws.Range("D214").Value = 12
$ws.Range("H214").Value = 1

$ws.Range("A215").Value = 'Islas Malvinas'
$ws.Range("D215").Value = 13
$ws.Range("H215").Value = 0

Write-Output "applied all changes"
